$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Completed varying queue size tests
#
# Two new queue-size columns ("1" and "10") are introduced into both
# summary tables on the sheet. For every data row, the existing values
# in columns E..H (old columns E,F,G,H) slide right by two columns to
# G..J, while the old C/D cell contents move into the newly vacated
# E/F slots, and C/D themselves receive the brand-new "1" / "10"
# queue-size data.
# -----------------------------------------------------------------

function Shift-Row([int]$row) {
    # Move existing data two columns to the right (rightmost first so
    # we never overwrite a source cell before it has been read).
    $ws.Range("H" + $row).Copy($ws.Range("J" + $row))
    $ws.Range("G" + $row).Copy($ws.Range("I" + $row))
    $ws.Range("F" + $row).Copy($ws.Range("H" + $row))
    $ws.Range("E" + $row).Copy($ws.Range("G" + $row))
    $ws.Range("D" + $row).Copy($ws.Range("F" + $row))
    $ws.Range("C" + $row).Copy($ws.Range("E" + $row))
}

# ----- Table 1 (rows 5-10) -----
Shift-Row 5
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 10

Shift-Row 6
$ws.Range("C6").Value = 1000
$ws.Range("D6").Value = 1000
# D6 gets a brand new visual style: same font/fill it already carried,
# plus a thin top border.
$ws.Range("D6").Borders.Item(8).LineStyle = 1

Shift-Row 7
$ws.Range("C7").Value = 8039
$ws.Range("D7").Value = 337

Shift-Row 8
$ws.Range("C8").Value = 0.05
$ws.Range("D8").Value = 16

Shift-Row 9
$ws.Range("C9").Value = 15426
$ws.Range("D9").Value = 49

# Row 10 only has a label in column C (plus blank filler cells) - move
# the existing "10 hours 24 mins" label two columns over to E10, give
# J10 an (empty) styled cell to match the rest of the table, drop the
# now-unused D10 cell, and write the new "4 hours 17 mins" label into C10.
$ws.Range("I10").Copy($ws.Range("J10"))
$ws.Range("C10").Copy($ws.Range("E10"))
$ws.Range("D10").Clear()
$ws.Range("C10").Value = "4 hours 17 mins"

# ----- Table 2 (rows 16-19) -----
Shift-Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 10

Shift-Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5

Shift-Row 18
$ws.Range("C18").Value = 2965285
$ws.Range("D18").Value = 3007046

Shift-Row 19
$ws.Range("C19").Value = 79
$ws.Range("D19").Value = 77

# Match the saved selection state from the authored workbook.
$ws.Range("C19").Select()
